$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "276.47"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.84%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.70"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.68%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.875"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.35%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06330"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.04%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.932"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.25%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.281"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "33.78%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8758"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.47%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1557"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "5.86%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05084"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.94%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07488"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.08%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02953"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-5.68%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09060"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.01%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001573"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.01%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006328"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.03%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006024"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "3.96%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.452"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.22%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.322"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.57%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.283"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.51%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.27%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1330"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.43%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.910"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.23%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04371"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.03%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001173"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.65%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004216"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.71%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.02%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001616"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-4.40%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04102"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.40%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006980"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.84%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1173"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.80%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002191"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.47%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01124"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-12.91%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005292"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.40%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-11.13%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.490"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-37.38%"
